$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "95.463.40"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.64%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.619.53"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.23%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "2.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +21.91%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "227.50"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.27%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "639.85"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.44%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.414"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.19%  "
$ws.Range("E9").Value = "  +2.69%  "
$ws.Range("E10").Value = "  -0.05%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "3.616.37"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.27%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "47.08"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +6.26%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.207"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.62%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000292"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.47%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.51"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.78%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.288.78"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.37%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "95.233.43"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.62%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.78"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.78%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "20.72"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +11.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.624.50"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.00"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.14%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.518"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.43%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "512.32"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.14%  "
$ws.Range("E24").Value = "  -4.94%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.246"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +26.37%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "119.95"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +17.63%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000204"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.95%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.76"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.41%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "12.68"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.45%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "12.80"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.35%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.92"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.88%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.03%  "
$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.80"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.42%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.180"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.42%  "
$ws.Range("B35").Value = "Binance-PegBSC-USD"
$ws.Range("C35").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.11%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "31.91"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.21%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.588"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.83%  "
$ws.Range("E38").Value = "  -0.02%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "599.83"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -7.25%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.36"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.35%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.83"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.38%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "40.59"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.18%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.159"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.04%  "
$ws.Range("E44").Value = "  +6.81%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0478"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.88%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.92"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.39%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.923"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.46%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.47"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.71%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.61"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.43%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.22"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.08%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "217.58"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.00%  "
